$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("E2").Value = 25.04000000000048
$ws.Range("G2").Value = 0.000000000000002664535259100376
$ws.Range("H2").Value = 0.00000000000004493137887894751
$ws.Range("K2").Value = 49.73486164912988
$ws.Range("L2").Value = "[37.10719197702518, 62.362531321234584]"
$ws.Range("M2").Value = 0.000000000000489608353859694
$ws.Range("N2").Value = 0.0000000000009792167077193881
$ws.Range("O2").Value = 1.37739497724958
$ws.Range("P2").Value = "[1.0880791372793475, 1.666710817219812]"
$ws.Range("S2").Value = 66.83793966144387
$ws.Range("T2").Value = "[59.10552611749495, 74.57035320539279]"
$ws.Range("W2").Value = 19.55075075075112
$ws.Range("X2").Value = 18.3977577577581
$ws.Range("Y2").Value = 20.70374374374414

# Row 3 updates
$ws.Range("E3").Value = 22
$ws.Range("G3").Value = 0.0000000006106859462562397
$ws.Range("H3").Value = 0.000000001623407324874787
$ws.Range("I3").Value = 0.4719533011521135
$ws.Range("K3").Value = 42.19988363933471
$ws.Range("L3").Value = "[27.741348845582635, 56.658418433086794]"
$ws.Range("M3").Value = 0.00000004826582644668065
$ws.Range("N3").Value = 0.00000004826582644668065
$ws.Range("O3").Value = -2.80510575275485
$ws.Range("P3").Value = "[-3.1950531892364666, -2.4151583162732337]"
$ws.Range("Q3").Value = 0
$ws.Range("R3").Value = 0
$ws.Range("S3").Value = 59.71576571621465
$ws.Range("T3").Value = "[51.1648730767763, 68.266658355653]"
$ws.Range("W3").Value = 9.821821821821821
$ws.Range("X3").Value = 8.456456456456454
$ws.Range("Y3").Value = 11.18718718718719
